# Update header labels on the LINE_TRIALS_URL sheet (also updates the
# bound table's column names automatically) and move the active cell
# selection to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LINE_TRIALS_URL")

$ws.Activate()

$ws.Range("B1").Value = "BOM_UNDER_TRIAL"
$ws.Range("G1").Value = "ORDER_NO."

$ws.Range("G1").Select()
